# Auto-generated edit script: updates market-data value cells (H:N) on several
# sheets per the scheduled-runner diff. Pure data values, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (item id 5470)
$ws.Range("H4").Value = 496.66666
$ws.Range("I4").Value = 345
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 345
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = -231
$ws.Range("N4").Value = -1028

# Row 46 (item id 4584)
$ws.Range("H46").Value = 2079.7693
$ws.Range("J46").Value = 1644.7
$ws.Range("L46").Value = 4934.1
$ws.Range("N46").Value = -5172.1

# Row 47 (item id 2169)
$ws.Range("H47").Value = 19500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 19500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 19500
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -21444

# Row 58 (item id 4606)
$ws.Range("H58").Value = 2094.4443
$ws.Range("I58").Value = 283.33334
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 850.0000200000001
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -700.0000200000001
$ws.Range("N58").Value = -9300

# Row 60 (item id 4584)
$ws.Range("H60").Value = 2079.7693
$ws.Range("J60").Value = 1644.7
$ws.Range("L60").Value = 4934.1
$ws.Range("N60").Value = -5902.1

# Row 62 (item id 27781)
$ws.Range("H62").Value = 4149
$ws.Range("I62").Value = 4266.6665
$ws.Range("K62").Value = 4266.6665
$ws.Range("M62").Value = -3642.6665

# Row 65 (item id 27781)
$ws.Range("H65").Value = 4149
$ws.Range("I65").Value = 4266.6665
$ws.Range("K65").Value = 21333.3325
$ws.Range("M65").Value = -18213.3325

# Row 98 (item id 36237)
$ws.Range("H98").Value = 1600.1111
$ws.Range("I98").Value = 1600.1111
$ws.Range("K98").Value = 1600.1111
$ws.Range("M98").Value = -102.1111000000001

# Row 122 (item id 36237)
$ws.Range("H122").Value = 1600.1111
$ws.Range("I122").Value = 1600.1111
$ws.Range("K122").Value = 4800.3333
$ws.Range("M122").Value = -2350.3333

# Row 129 (item id 36115)
$ws.Range("H129").Value = 882.60785
$ws.Range("J129").Value = 946
$ws.Range("L129").Value = 2838
$ws.Range("N129").Value = -12838

# Row 135 (item id 44047)
$ws.Range("H135").Value = 55556856
$ws.Range("I135").Value = 20834460
$ws.Range("J135").Value = 333336000
$ws.Range("K135").Value = 187510140
$ws.Range("L135").Value = 3000024000
$ws.Range("M135").Value = -187507605
$ws.Range("N135").Value = -3000029070


$ws = $wb.Worksheets.Item("ARM")
# Row 5 (item id 5091)
$ws.Range("H5").Value = 40
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 40
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -264

# Row 56 (item id 2504)
$ws.Range("H56").Value = 21946.25
$ws.Range("J56").Value = 21946.25
$ws.Range("L56").Value = 21946.25
$ws.Range("N56").Value = -23430.25

# Row 74 (item id 44000)
$ws.Range("H74").Value = 1807.9615
$ws.Range("I74").Value = 1458.625
$ws.Range("K74").Value = 1458.625
$ws.Range("M74").Value = -584.625

# Row 77 (item id 44000)
$ws.Range("H77").Value = 1807.9615
$ws.Range("I77").Value = 1458.625
$ws.Range("K77").Value = 7293.125
$ws.Range("M77").Value = -2925.125


$ws = $wb.Worksheets.Item("BSM")
# Row 4 (item id 5091)
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 40
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -270

# Row 22 (item id 5092)
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127

# Row 105 (item id 19947)
$ws.Range("H105").Value = 3433.6604
$ws.Range("I105").Value = 2984.625
$ws.Range("J105").Value = 4815.3076
$ws.Range("K105").Value = 2984.625
$ws.Range("L105").Value = 4815.3076
$ws.Range("M105").Value = -1237.625
$ws.Range("N105").Value = -8309.3076

# Row 134 (item id 43998)
$ws.Range("H134").Value = 5940.3228
$ws.Range("I134").Value = 8399
$ws.Range("K134").Value = 25197
$ws.Range("M134").Value = -22662


$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item id 5367)
$ws.Range("H22").Value = 105.90909
$ws.Range("I22").Value = 73.57143000000001
$ws.Range("J22").Value = 162.5
$ws.Range("K22").Value = 73.57143000000001
$ws.Range("L22").Value = 162.5
$ws.Range("M22").Value = 276.42857
$ws.Range("N22").Value = -862.5

# Row 58 (item id 44021)
$ws.Range("H58").Value = 3956157.5
$ws.Range("I58").Value = 10103510
$ws.Range("J58").Value = 4288.0713
$ws.Range("K58").Value = 10103510
$ws.Range("L58").Value = 4288.0713
$ws.Range("M58").Value = -10103307
$ws.Range("N58").Value = -4694.0713

# Row 107 (item id 27689)
$ws.Range("H107").Value = 1071.3636
$ws.Range("I107").Value = 1120.6875
$ws.Range("J107").Value = 939.8333
$ws.Range("K107").Value = 1120.6875
$ws.Range("L107").Value = 939.8333
$ws.Range("M107").Value = 799.3125
$ws.Range("N107").Value = -4779.8333

# Row 132 (item id 44019)
$ws.Range("H132").Value = 2471.9
$ws.Range("I132").Value = 1998.7693
$ws.Range("J132").Value = 3350.5715
$ws.Range("K132").Value = 5996.3079
$ws.Range("L132").Value = 10051.7145
$ws.Range("M132").Value = -3466.3079
$ws.Range("N132").Value = -15111.7145

# Row 134 (item id 44020)
$ws.Range("H134").Value = 3464.2258
$ws.Range("I134").Value = 2473.9375
$ws.Range("J134").Value = 4520.533
$ws.Range("K134").Value = 7421.8125
$ws.Range("L134").Value = 13561.599
$ws.Range("M134").Value = -4886.8125
$ws.Range("N134").Value = -18631.599

# Row 136 (item id 44021)
$ws.Range("H136").Value = 3956157.5
$ws.Range("I136").Value = 10103510
$ws.Range("J136").Value = 4288.0713
$ws.Range("K136").Value = 30310530
$ws.Range("L136").Value = 12864.2139
$ws.Range("M136").Value = -30307980
$ws.Range("N136").Value = -17964.2139


$ws = $wb.Worksheets.Item("GSM")
# Row 102 (item id 36169)
$ws.Range("H102").Value = 3665.6
$ws.Range("I102").Value = 2854.72
$ws.Range("K102").Value = 2854.72
$ws.Range("M102").Value = -1232.72


$ws = $wb.Worksheets.Item("LTW")
# Row 46 (item id 5282)
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 68 (item id 12563)
$ws.Range("H68").Value = 1250
$ws.Range("I68").Value = 1250
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1250
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -501
$ws.Range("N68").ClearContents()

# Row 71 (item id 12563)
$ws.Range("H71").Value = 1250
$ws.Range("I71").Value = 1250
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6250
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2506
$ws.Range("N71").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
# Row 62 (item id 12589)
$ws.Range("H62").Value = 4125
$ws.Range("I62").Value = 3960
$ws.Range("J62").Value = 4400
$ws.Range("K62").Value = 3960
$ws.Range("L62").Value = 4400
$ws.Range("M62").Value = -3336
$ws.Range("N62").Value = -5648

# Row 65 (item id 12589)
$ws.Range("H65").Value = 4125
$ws.Range("I65").Value = 3960
$ws.Range("J65").Value = 4400
$ws.Range("K65").Value = 19800
$ws.Range("L65").Value = 22000
$ws.Range("M65").Value = -16680
$ws.Range("N65").Value = -28240

# Row 122 (item id 36208)
$ws.Range("H122").Value = 2106.7693
$ws.Range("I122").Value = 1896.4445
$ws.Range("J122").Value = 2580
$ws.Range("K122").Value = 5689.333500000001
$ws.Range("L122").Value = 7740
$ws.Range("M122").Value = -3239.333500000001
$ws.Range("N122").Value = -12640

# Row 132 (item id 44029)
$ws.Range("H132").Value = 1259.7291
$ws.Range("I132").Value = 640.0909
$ws.Range("J132").Value = 2622.9333
$ws.Range("K132").Value = 1920.2727
$ws.Range("L132").Value = 7868.7999
$ws.Range("M132").Value = 609.7273
$ws.Range("N132").Value = -12928.7999

